$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: In the paragraph about the terms of the contract, add a new
# sentence about menial tasks, right before "Each human is given a choice..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "depending on the legend specifying the terms. Each human is given a choice",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "depending on the legend specifying the terms. Terms can also include menial tasks, like picking up trash, if such a legend desires it (i.e. Gods of nature for instance). Each human is given a choice",
    2) | Out-Null

# ---------------------------------------------------------------------------
# Edit 2: Expand on what "powers and weaknesses" means, rewriting the
# sentence that used to read:
#   " In other words, they gain both powers and weaknesses based on that
#     legend's background."
# into a longer explanation of special abilities/weaknesses.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    " In other words, they gain both powers and weaknesses based on that legend's background.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " In other words, they gain both special abilities and weaknesses based on that legend contracting them. Special abilities can range from materializing special weapons, physical enhancements, magical abilities, or any form of supernatural power related to the legend. The weaknesses gained are also based from the contracted legend and could range from simple exposure to rain or being damaged by Holy Light magic, both of which can be fatal if it was originally a weakness/cause of death from the original legend. This will signify the importance of keeping their contracted legend a secret, lest their weakness will be exposed.",
    2) | Out-Null

# ---------------------------------------------------------------------------
# Edit 3: Split the paragraph into two right before "However, overuse of",
# turning that sentence into its own paragraph, and extend it with
# "regardless of the legend,". The leading space that used to separate it
# from the previous sentence is consumed by the new paragraph mark.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    " However, overuse of",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "However, regardless of the legend, overuse of",
    2) | Out-Null

$splitRng = $d.Content
$splitRng.Find.Execute("However, regardless of the legend, overuse of", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPoint = $d.Range($splitRng.Start, $splitRng.Start)
$splitPoint.InsertParagraphBefore()
